# Update path-node data (columns A-D, rows 4-75) to reflect the
# recomputed "accelpressure" path-finding results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 3).Value = 14
$ws.Cells.Item(4, 4).Value = 4.0404040404040407
# Row 5
$ws.Cells.Item(5, 2).Value = 670
$ws.Cells.Item(5, 3).Value = 16
$ws.Cells.Item(5, 4).Value = 4.7070707070707076
# Row 6
$ws.Cells.Item(6, 2).Value = 680
$ws.Cells.Item(6, 3).Value = 18
$ws.Cells.Item(6, 4).Value = 5.2953060011883544
# Row 7
$ws.Cells.Item(7, 2).Value = 690
$ws.Cells.Item(7, 3).Value = 20
$ws.Cells.Item(7, 4).Value = 5.8216217906620384
# Row 8
$ws.Cells.Item(8, 1).Value = 890
$ws.Cells.Item(8, 2).Value = 710
$ws.Cells.Item(8, 3).Value = 23
$ws.Cells.Item(8, 4).Value = 6.7518543488015732
# Row 9
$ws.Cells.Item(9, 1).Value = 890
$ws.Cells.Item(9, 2).Value = 820
$ws.Cells.Item(9, 3).Value = 34
$ws.Cells.Item(9, 4).Value = 10.611503471608589
# Row 10
$ws.Cells.Item(10, 1).Value = 890
$ws.Cells.Item(10, 2).Value = 830
$ws.Cells.Item(10, 3).Value = 35
$ws.Cells.Item(10, 4).Value = 10.901358544072361
# Row 11
$ws.Cells.Item(11, 1).Value = 880
$ws.Cells.Item(11, 2).Value = 930
$ws.Cells.Item(11, 3).Value = 26
$ws.Cells.Item(11, 4).Value = 14.19639973132511
# Row 12
$ws.Cells.Item(12, 1).Value = 860
$ws.Cells.Item(12, 2).Value = 980
$ws.Cells.Item(12, 3).Value = 25
$ws.Cells.Item(12, 4).Value = 16.308229067456288
# Row 13
$ws.Cells.Item(13, 1).Value = 850
$ws.Cells.Item(13, 2).Value = 990
$ws.Cells.Item(13, 3).Value = 19
$ws.Cells.Item(13, 4).Value = 16.951053413989509
# Row 14
$ws.Cells.Item(14, 1).Value = 840
$ws.Cells.Item(14, 2).Value = 1000
$ws.Cells.Item(14, 3).Value = 21
$ws.Cells.Item(14, 4).Value = 17.658160195176059
# Row 15
$ws.Cells.Item(15, 1).Value = 820
$ws.Cells.Item(15, 2).Value = 1000
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = 19.476342013357879
# Row 16
$ws.Cells.Item(16, 1).Value = 800
$ws.Cells.Item(16, 2).Value = 990
$ws.Cells.Item(16, 3).Value = 12
$ws.Cells.Item(16, 4).Value = 22.916446594126779
# Row 17
$ws.Cells.Item(17, 1).Value = 790
$ws.Cells.Item(17, 2).Value = 980
$ws.Cells.Item(17, 3).Value = 15
$ws.Cells.Item(17, 4).Value = 23.964012195884631
# Row 18
$ws.Cells.Item(18, 1).Value = 780
$ws.Cells.Item(18, 2).Value = 970
$ws.Cells.Item(18, 3).Value = 18
$ws.Cells.Item(18, 4).Value = 24.821111324595599
# Row 19
$ws.Cells.Item(19, 1).Value = 760
$ws.Cells.Item(19, 2).Value = 930
$ws.Cells.Item(19, 3).Value = 23
$ws.Cells.Item(19, 4).Value = 27.00264105874173
# Row 20
$ws.Cells.Item(20, 1).Value = 730
$ws.Cells.Item(20, 2).Value = 870
$ws.Cells.Item(20, 3).Value = 26
$ws.Cells.Item(20, 4).Value = 29.740683480170048
# Row 21
$ws.Cells.Item(21, 1).Value = 710
$ws.Cells.Item(21, 2).Value = 850
$ws.Cells.Item(21, 3).Value = 23
$ws.Cells.Item(21, 4).Value = 30.895143531086859
# Row 22
$ws.Cells.Item(22, 1).Value = 630
$ws.Cells.Item(22, 2).Value = 820
$ws.Cells.Item(22, 3).Value = 19.27204019518533
$ws.Cells.Item(22, 4).Value = 34.937533586483653
# Row 23
$ws.Cells.Item(23, 1).Value = 620
$ws.Cells.Item(23, 3).Value = 20.379030654534951
$ws.Cells.Item(23, 4).Value = 35.441933583038633
# Row 24
$ws.Cells.Item(24, 1).Value = 590
$ws.Cells.Item(24, 2).Value = 830
$ws.Cells.Item(24, 3).Value = 24
$ws.Cells.Item(24, 4).Value = 36.867056039986757
# Row 25
$ws.Cells.Item(25, 1).Value = 520
$ws.Cells.Item(25, 2).Value = 880
$ws.Cells.Item(25, 3).Value = 25.79757424973122
$ws.Cells.Item(25, 4).Value = 40.321973431706517
# Row 26
$ws.Cells.Item(26, 1).Value = 490
$ws.Cells.Item(26, 2).Value = 890
$ws.Cells.Item(26, 3).Value = 23.529934170619271
$ws.Cells.Item(26, 4).Value = 41.604129275954989
# Row 27
$ws.Cells.Item(27, 1).Value = 460
$ws.Cells.Item(27, 2).Value = 890
$ws.Cells.Item(27, 3).Value = 24
$ws.Cells.Item(27, 4).Value = 42.866491638684479
# Row 28
$ws.Cells.Item(28, 1).Value = 400
$ws.Cells.Item(28, 2).Value = 860
$ws.Cells.Item(28, 3).Value = 17
$ws.Cells.Item(28, 4).Value = 46.138786239903681
# Row 29
$ws.Cells.Item(29, 1).Value = 390
$ws.Cells.Item(29, 2).Value = 850
$ws.Cells.Item(29, 3).Value = 19
$ws.Cells.Item(29, 4).Value = 46.924460441222067
# Row 30
$ws.Cells.Item(30, 1).Value = 370
$ws.Cells.Item(30, 2).Value = 800
$ws.Cells.Item(30, 3).Value = 20.379030654534951
$ws.Cells.Item(30, 4).Value = 49.659502273194811
# Row 31
$ws.Cells.Item(31, 1).Value = 370
$ws.Cells.Item(31, 2).Value = 770
$ws.Cells.Item(31, 3).Value = 21
$ws.Cells.Item(31, 4).Value = 51.109512074075063
# Row 32
$ws.Cells.Item(32, 1).Value = 380
$ws.Cells.Item(32, 2).Value = 740
$ws.Cells.Item(32, 3).Value = 24
$ws.Cells.Item(32, 4).Value = 52.514968811927673
# Row 33
$ws.Cells.Item(33, 1).Value = 420
$ws.Cells.Item(33, 2).Value = 680
$ws.Cells.Item(33, 3).Value = 27
$ws.Cells.Item(33, 4).Value = 55.342852165232763
# Row 34
$ws.Cells.Item(34, 1).Value = 490
$ws.Cells.Item(34, 2).Value = 600
$ws.Cells.Item(34, 3).Value = 30
$ws.Cells.Item(34, 4).Value = 59.07272788899931
# Row 35
$ws.Cells.Item(35, 1).Value = 540
$ws.Cells.Item(35, 2).Value = 520
$ws.Cells.Item(35, 3).Value = 25
$ws.Cells.Item(35, 4).Value = 62.503266482474437
# Row 36
$ws.Cells.Item(36, 1).Value = 550
$ws.Cells.Item(36, 2).Value = 460
$ws.Cells.Item(36, 3).Value = 20.961723555504101
$ws.Cells.Item(36, 4).Value = 65.150148304890251
# Row 37
$ws.Cells.Item(37, 1).Value = 550
$ws.Cells.Item(37, 2).Value = 430
$ws.Cells.Item(37, 3).Value = 24
$ws.Cells.Item(37, 4).Value = 66.484616720584583
# Row 38
$ws.Cells.Item(38, 1).Value = 540
$ws.Cells.Item(38, 2).Value = 410
$ws.Cells.Item(38, 3).Value = 17
$ws.Cells.Item(38, 4).Value = 67.575381587657645
# Row 39
$ws.Cells.Item(39, 1).Value = 530
$ws.Cells.Item(39, 2).Value = 400
$ws.Cells.Item(39, 3).Value = 19
$ws.Cells.Item(39, 4).Value = 68.361055788976032
# Row 40
$ws.Cells.Item(40, 1).Value = 500
$ws.Cells.Item(40, 2).Value = 380
$ws.Cells.Item(40, 3).Value = 24
$ws.Cells.Item(40, 4).Value = 70.038056382215103
# Row 41
$ws.Cells.Item(41, 1).Value = 460
$ws.Cells.Item(41, 2).Value = 370
$ws.Cells.Item(41, 3).Value = 21.62096934546506
$ws.Cells.Item(41, 4).Value = 71.845604834066449
# Row 42
$ws.Cells.Item(42, 1).Value = 450
$ws.Cells.Item(42, 2).Value = 370
$ws.Cells.Item(42, 3).Value = 23
$ws.Cells.Item(42, 4).Value = 72.293824590234351
# Row 43
$ws.Cells.Item(43, 1).Value = 400
$ws.Cells.Item(43, 2).Value = 380
$ws.Cells.Item(43, 3).Value = 27
$ws.Cells.Item(43, 4).Value = 74.333432395671466
# Row 44
$ws.Cells.Item(44, 1).Value = 370
$ws.Cells.Item(44, 2).Value = 390
$ws.Cells.Item(44, 3).Value = 24
$ws.Cells.Item(44, 4).Value = 75.573541282012002
# Row 45
$ws.Cells.Item(45, 1).Value = 350
$ws.Cells.Item(45, 2).Value = 390
$ws.Cells.Item(45, 3).Value = 24
$ws.Cells.Item(45, 4).Value = 76.40687461534533
# Row 46
$ws.Cells.Item(46, 1).Value = 330
$ws.Cells.Item(46, 2).Value = 380
$ws.Cells.Item(46, 3).Value = 17
$ws.Cells.Item(46, 4).Value = 77.497639482418393
# Row 47
$ws.Cells.Item(47, 1).Value = 300
$ws.Cells.Item(47, 2).Value = 350
$ws.Cells.Item(47, 3).Value = 20
$ws.Cells.Item(47, 4).Value = 79.790958772753143
# Row 48
$ws.Cells.Item(48, 1).Value = 280
$ws.Cells.Item(48, 2).Value = 300
$ws.Cells.Item(48, 3).Value = 20
$ws.Cells.Item(48, 4).Value = 82.483541176320401
# Row 49
$ws.Cells.Item(49, 1).Value = 280
$ws.Cells.Item(49, 2).Value = 290
$ws.Cells.Item(49, 3).Value = 21
$ws.Cells.Item(49, 4).Value = 82.971346054369178
# Row 50
$ws.Cells.Item(50, 1).Value = 290
$ws.Cells.Item(50, 2).Value = 260
$ws.Cells.Item(50, 3).Value = 24.61337882406934
$ws.Cells.Item(50, 4).Value = 84.357903131127415
# Row 51
$ws.Cells.Item(51, 1).Value = 320
$ws.Cells.Item(51, 2).Value = 220
$ws.Cells.Item(51, 3).Value = 24
$ws.Cells.Item(51, 4).Value = 86.414950026836877
# Row 52
$ws.Cells.Item(52, 1).Value = 340
$ws.Cells.Item(52, 2).Value = 200
$ws.Cells.Item(52, 3).Value = 27
$ws.Cells.Item(52, 4).Value = 87.524137134580485
# Row 53
$ws.Cells.Item(53, 1).Value = 350
$ws.Cells.Item(53, 2).Value = 190
$ws.Cells.Item(53, 3).Value = 29
$ws.Cells.Item(53, 4).Value = 88.029213406856584
# Row 54
$ws.Cells.Item(54, 1).Value = 430
$ws.Cells.Item(54, 2).Value = 150
$ws.Cells.Item(54, 3).Value = 24.61337882406934
$ws.Cells.Item(54, 4).Value = 91.365795471325796
# Row 55
$ws.Cells.Item(55, 1).Value = 480
$ws.Cells.Item(55, 2).Value = 140
$ws.Cells.Item(55, 3).Value = 27
$ws.Cells.Item(55, 4).Value = 93.341647325928619
# Row 56
$ws.Cells.Item(56, 1).Value = 490
$ws.Cells.Item(56, 2).Value = 140
$ws.Cells.Item(56, 3).Value = 28
$ws.Cells.Item(56, 4).Value = 93.705283689564979
# Row 57
$ws.Cells.Item(57, 1).Value = 530
$ws.Cells.Item(57, 2).Value = 140
$ws.Cells.Item(57, 3).Value = 32
$ws.Cells.Item(57, 4).Value = 95.038617022898308
# Row 58
$ws.Cells.Item(58, 1).Value = 540
$ws.Cells.Item(58, 2).Value = 140
$ws.Cells.Item(58, 3).Value = 33
$ws.Cells.Item(58, 4).Value = 95.346309330590614
# Row 59
$ws.Cells.Item(59, 1).Value = 560
$ws.Cells.Item(59, 2).Value = 140
$ws.Cells.Item(59, 3).Value = 35
$ws.Cells.Item(59, 4).Value = 95.934544624708266
# Row 60
$ws.Cells.Item(60, 1).Value = 580
$ws.Cells.Item(60, 2).Value = 140
$ws.Cells.Item(60, 3).Value = 37
$ws.Cells.Item(60, 4).Value = 96.490100180263823
# Row 61
$ws.Cells.Item(61, 1).Value = 590
$ws.Cells.Item(61, 2).Value = 140
$ws.Cells.Item(61, 3).Value = 38
$ws.Cells.Item(61, 4).Value = 96.756766846930489
# Row 62
$ws.Cells.Item(62, 1).Value = 650
$ws.Cells.Item(62, 2).Value = 140
$ws.Cells.Item(62, 3).Value = 37
$ws.Cells.Item(62, 4).Value = 98.356766846930483
# Row 63
$ws.Cells.Item(63, 1).Value = 710
$ws.Cells.Item(63, 2).Value = 150
$ws.Cells.Item(63, 3).Value = 32.589046926199423
$ws.Cells.Item(63, 4).Value = 100.1049622308245
# Row 64
$ws.Cells.Item(64, 1).Value = 810
$ws.Cells.Item(64, 2).Value = 200
$ws.Cells.Item(64, 3).Value = 24
$ws.Cells.Item(64, 4).Value = 104.0563770344983
# Row 65
$ws.Cells.Item(65, 1).Value = 870
$ws.Cells.Item(65, 2).Value = 260
$ws.Cells.Item(65, 3).Value = 24.61337882406934
$ws.Cells.Item(65, 4).Value = 107.54730137651811
# Row 66
$ws.Cells.Item(66, 1).Value = 910
$ws.Cells.Item(66, 2).Value = 330
$ws.Cells.Item(66, 3).Value = 27.439894486170861
$ws.Cells.Item(66, 4).Value = 110.6449961939975
# Row 67
$ws.Cells.Item(67, 1).Value = 920
$ws.Cells.Item(67, 2).Value = 370
$ws.Cells.Item(67, 3).Value = 27
$ws.Cells.Item(67, 4).Value = 112.1597330113751
# Row 68
$ws.Cells.Item(68, 1).Value = 920
$ws.Cells.Item(68, 2).Value = 380
$ws.Cells.Item(68, 3).Value = 28
$ws.Cells.Item(68, 4).Value = 112.5233693750114
# Row 69
$ws.Cells.Item(69, 1).Value = 920
$ws.Cells.Item(69, 2).Value = 420
$ws.Cells.Item(69, 3).Value = 32
$ws.Cells.Item(69, 4).Value = 113.8567027083448
# Row 70
$ws.Cells.Item(70, 1).Value = 920
$ws.Cells.Item(70, 2).Value = 430
$ws.Cells.Item(70, 3).Value = 33
$ws.Cells.Item(70, 4).Value = 114.1643950160371
# Row 71
$ws.Cells.Item(71, 1).Value = 920
$ws.Cells.Item(71, 2).Value = 440
$ws.Cells.Item(71, 3).Value = 34
$ws.Cells.Item(71, 4).Value = 114.46290247872361
# Row 72
$ws.Cells.Item(72, 1).Value = 920
$ws.Cells.Item(72, 2).Value = 450
$ws.Cells.Item(72, 3).Value = 35
$ws.Cells.Item(72, 4).Value = 114.7527575511874
# Row 73
$ws.Cells.Item(73, 1).Value = 920
$ws.Cells.Item(73, 2).Value = 500
$ws.Cells.Item(73, 3).Value = 39
$ws.Cells.Item(73, 4).Value = 116.1041089025388
# Row 74
$ws.Cells.Item(74, 1).Value = 920
$ws.Cells.Item(74, 2).Value = 520
$ws.Cells.Item(74, 3).Value = 41
$ws.Cells.Item(74, 4).Value = 116.6041089025388
# Row 75
$ws.Cells.Item(75, 1).Value = 920
$ws.Cells.Item(75, 2).Value = 530
$ws.Cells.Item(75, 3).Value = 41
$ws.Cells.Item(75, 4).Value = 116.84801134156319
